$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 25.98413200434529
$ws.Range("B2").Value = 674
$ws.Range("C2").Value = 568
$ws.Range("D2").Value = 800
$ws.Range("E2").Value = 448
$ws.Range("F2").Value = "nao_busque"

# Row 3
$ws.Range("A3").Value = 19.95105704378432
$ws.Range("B3").Value = 907
$ws.Range("C3").Value = 528
$ws.Range("D3").Value = 800
$ws.Range("E3").Value = 448
$ws.Range("F3").Value = "nao_busque"

# Row 4
$ws.Range("A4").Value = 122.6395717848804
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 642
$ws.Range("D4").Value = 800
$ws.Range("E4").Value = 448
$ws.Range("F4").Value = "busque"

# Row 5
$ws.Range("A5").Value = 123.7064184900456
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 663
$ws.Range("D5").Value = 800
$ws.Range("E5").Value = 448
$ws.Range("F5").Value = "busque"

# Row 6
$ws.Range("A6").Value = 30.37104468118742
$ws.Range("B6").Value = 1001
$ws.Range("C6").Value = 479
$ws.Range("D6").Value = 800
$ws.Range("E6").Value = 448
$ws.Range("F6").Value = "nao_busque"

# Row 7
$ws.Range("A7").Value = 26.96379937672894
$ws.Range("B7").Value = 649
$ws.Range("C7").Value = 547
$ws.Range("D7").Value = 800
$ws.Range("E7").Value = 448
$ws.Range("F7").Value = "nao_busque"

# Row 8
$ws.Range("A8").Value = 70.46607276494409
$ws.Range("B8").Value = 338
$ws.Range("C8").Value = 544
$ws.Range("D8").Value = 800
$ws.Range("E8").Value = 448
$ws.Range("F8").Value = "busque"

# Row 9
$ws.Range("A9").Value = 51.97877640881225
$ws.Range("B9").Value = 452
$ws.Range("C9").Value = 441
$ws.Range("D9").Value = 800
$ws.Range("E9").Value = 448
$ws.Range("F9").Value = "busque"

# Row 10
$ws.Range("A10").Value = 119.530350585183
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 474
$ws.Range("D10").Value = 800
$ws.Range("E10").Value = 448
$ws.Range("F10").Value = "busque"

# Row 11
$ws.Range("A11").Value = 128.4761038825038
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 767
$ws.Range("D11").Value = 800
$ws.Range("E11").Value = 448
$ws.Range("F11").Value = "busque"
